$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B24: was text "3" typed as an inline string -> should become a real number 3 ---
$ws.Range("B24").Value = 3

# --- Add new row 25 ---
$ws.Range("A25").Value = "Ying Tang"

# B25 must stay TEXT "2" (not be auto-converted to the number 2, and without
# picking up a new NumberFormat/quote-prefix style). Build it as a formula
# that evaluates to the text string "2", then paste-special just the value
# over the target cell - this preserves the Text cell type without
# introducing any new style record.
$helper = $ws.Range("ZZ1")
$helper.Formula = "=""2"""
$helper.Copy() | Out-Null
$ws.Range("B25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$helper.Clear() | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C25").Value = "Cons,it does not seem to involve,it clearly does not fit "
$ws.Range("D25").Value = "CRT"
$ws.Range("E25").Value = "MET"
$ws.Range("F25").Value = "b49eb73e-9ff0-45de-a177-7d78dc315c92"
$ws.Range("G25").Value = "2rHk2kZ5knTJ6_annotated.xlsx"
$ws.Range("H25").Value = "Cons: - it does not seem to involve any learning, it clearly does not fit at ICLR."
